$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to
# Text format first, otherwise Excel auto-converts the assigned string
# into a numeric value (losing the exact text representation).
$ws.Range("D2").Value = '37.047.03'
$ws.Range("E2").Value = '  -1.61%  '
$ws.Range("D3").Value = '2.022.05'
$ws.Range("E3").Value = '  -3.15%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.39'
$ws.Range("E5").Value = '  -3.24%  '
$ws.Range("E6").Value = '  -4.58%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.88'
$ws.Range("E8").Value = '  -5.66%  '
$ws.Range("E9").Value = '  -2.95%  '
$ws.Range("E10").Value = '  +1.38%  '
$ws.Range("E11").Value = '  -3.84%  '
$ws.Range("D12").Value = '2.323.25'
$ws.Range("E12").Value = '  -2.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.27'
$ws.Range("E13").Value = '  -6.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.56'
$ws.Range("E14").Value = '  -2.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.744'
$ws.Range("E15").Value = '  -3.94%  '
$ws.Range("E16").Value = '  -4.09%  '
$ws.Range("D17").Value = '2.020.05'
$ws.Range("E17").Value = '  -3.05%  '
$ws.Range("D18").Value = '36.983.31'
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.85'
$ws.Range("E20").Value = '  -2.87%  '
$ws.Range("D21").Value = '0.0₃0828'
$ws.Range("E21").Value = '  -0.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.68'
$ws.Range("E22").Value = '  -1.68%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +1.95%  '
$ws.Range("E25").Value = '  -5.18%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.33'
$ws.Range("E26").Value = '  -3.80%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.48'
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.127'
$ws.Range("E28").Value = '  -4.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.77'
$ws.Range("E29").Value = '  -3.87%  '
$ws.Range("E30").Value = '  -3.70%  '
$ws.Range("E31").Value = '  -5.13%  '
$ws.Range("E32").Value = '  -4.17%  '
$ws.Range("E33").Value = '  -4.54%  '
$ws.Range("E34").Value = '  -4.92%  '
$ws.Range("E35").Value = '  -5.19%  '
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.18'
$ws.Range("E38").Value = '  -4.90%  '
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("D40").Value = '1.502.89'
$ws.Range("E40").Value = '  +2.63%  '
$ws.Range("E41").Value = '  -6.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.08'
$ws.Range("E42").Value = '  +1.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '95.08'
$ws.Range("E43").Value = '  -6.10%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0927'
$ws.Range("E44").Value = '  -3.60%  '
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  -3.17%  '
$ws.Range("E46").Value = '  -6.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.24'
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("E48").Value = '  -4.81%  '
$ws.Range("E49").Value = '  -1.82%  '
$ws.Range("B50").Value = 'FTXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.77'
$ws.Range("E50").Value = '  -6.24%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.211.19'
$ws.Range("E51").Value = '  -2.77%  '
